$wb = $excel.ActiveWorkbook

# --- Sheet "ev_charging_uc": update the two timeslice-grouping strings (C13/C14) ---
# G7 (=C14) and G8 (=C13) recalc automatically once the source cells change.
$wsEv = $wb.Worksheets.Item("ev_charging_uc")
$wsEv.Range("C13").Value = "FaP,SaP,WaP,SaD,FaD,RaP,RaD,WaD"
$wsEv.Range("C14").Value = "RaN,FaP,SaP,RaP,FaN,WaP,SaN,WaN"

# --- Sheet "re_profiles": rows 4-7 (columns M:O) rotate up by one row ---
$wsRe = $wb.Worksheets.Item("re_profiles")
$wsRe.Range("M4").Value = "R"
$wsRe.Range("N4").Value = 0.34481908618716439
$wsRe.Range("M5").Value = "S"
$wsRe.Range("N5").Value = 0.27551721102209703
$wsRe.Range("M6").Value = "F"
$wsRe.Range("N6").Value = 0.28270094198432955
$wsRe.Range("M7").Value = "W"
$wsRe.Range("N7").Value = 0.29696276080640904
